# Add a new sheet "2022-Q1" positioned between "2021-Q2" and "总计",
# populate it with fund-holding data, and insert a new summary row at the
# top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# --- 1. Insert new worksheet right after "2021-Q2" ------------------------
$ws = $sheets.Add($null, $sheets.Item("2021-Q2"))
$ws.Name = "2022-Q1"

# NOTE: fetch worksheet references *after* the sheet insertion above -
# inserting a sheet shifts worksheet positions, so references grabbed
# beforehand can end up pointing at the wrong sheet.
$ws = $sheets.Item("2022-Q1")
$tot = $sheets.Item("总计")

# --- 2. Copy the header/index-column look-and-feel from "总计" ------------
# (keeps the same bold/bordered/centred style used by the existing sheets)
$tot.Range("B1:D1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$tot.Range("A2").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Header row ----------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- 4. Data rows ------------------------------------------------------------
# Fund codes (B) and the decimal-figure columns (D:G) are stored as plain
# text in the source data (e.g. leading zeros in codes, trailing zeros in
# ratios), so force a text number format before writing the values.
$ws.Range("B2:B13").NumberFormat = "@"
$ws.Range("D2:G13").NumberFormat = "@"

$data = @(
    @("040007", "华安中小盘成长混合",               "22.83", "83.32", "3.24", "0.7397", 3),
    @("040025", "华安科技动力混合",                 "10.15", "82.51", "3.03", "0.3075", 5),
    @("001072", "华安智能装备主题股票",             "6.91",  "85.85", "2.61", "0.1804", 9),
    @("006122", "华安低碳生活混合",                 "4.43",  "86.68", "3.25", "0.1440", 6),
    @("011550", "湘财创新成长一年持有期混合A",       "2.62",  "93.51", "3.87", "0.1014", 10),
    @("008635", "华安科技创新混合",                 "3.29",  "83.14", "3.01", "0.0990", 3),
    @("001569", "泰信国策驱动灵活配置混合",         "2.19",  "85.17", "3.17", "0.0694", 10),
    @("003835", "鹏华沪深港新兴成长灵活配置混合",   "0.61",  "82.70", "4.64", "0.0283", 3),
    @("010076", "湘财长弘灵活配置混合A",             "0.32",  "91.72", "3.56", "0.0114", 10),
    @("011551", "湘财创新成长一年持有期混合C",       "0.28",  "93.51", "3.87", "0.0108", 10),
    @("010077", "湘财长弘灵活配置混合C",             "0.10",  "91.72", "3.56", "0.0036", 10),
    @("003739", "新华鑫弘灵活配置混合",             "0.01",  "42.81", "0.47", "0.0000", 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $fund = $data[$i]

    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $fund[0]
    $ws.Cells.Item($row, 3).Value = $fund[1]
    $ws.Cells.Item($row, 4).Value = $fund[2]
    $ws.Cells.Item($row, 5).Value = $fund[3]
    $ws.Cells.Item($row, 6).Value = $fund[4]
    $ws.Cells.Item($row, 7).Value = $fund[5]
    $ws.Cells.Item($row, 8).Value = $fund[6]
}

# --- 5. Insert a new top data row into "总计" for the 2022-Q1 summary ------
$tot.Rows.Item(2).Insert()
# Excel's row-insert inherits the row-above's formatting; the new data row
# should be plain (unstyled), matching the other data rows in this sheet.
$tot.Rows.Item(2).ClearFormats()

$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 12
$tot.Cells.Item(2, 4).Value = 1.7

# Column A keeps the same centred/bold index-column style as the rest of
# the sheet (the pre-existing rows already carry it).
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$tot.Cells.Item(3, 1).Value = 1

# --- 6. Restore the originally active sheet/tab ----------------------------
$sheets.Item("2021-Q2").Activate()
